$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (69 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1697183.8
$ws.Range("J17").Value = 1726434.1
$ws.Range("L17").Value = 5179302.300000001
$ws.Range("N17").Value = -5179638.300000001
$ws.Range("H62").Value = 2908.611
$ws.Range("I62").Value = 2409.5454
$ws.Range("J62").Value = 3692.8572
$ws.Range("K62").Value = 2409.5454
$ws.Range("L62").Value = 3692.8572
$ws.Range("M62").Value = -1785.5454
$ws.Range("N62").Value = -4940.8572
$ws.Range("H65").Value = 2908.611
$ws.Range("I65").Value = 2409.5454
$ws.Range("J65").Value = 3692.8572
$ws.Range("K65").Value = 12047.727
$ws.Range("L65").Value = 18464.286
$ws.Range("M65").Value = -8927.726999999999
$ws.Range("N65").Value = -24704.286
$ws.Range("H74").Value = 3885.8333
$ws.Range("I74").Value = 2377.5
$ws.Range("J74").Value = 4074.375
$ws.Range("K74").Value = 2377.5
$ws.Range("L74").Value = 4074.375
$ws.Range("M74").Value = -1441.5
$ws.Range("N74").Value = -5946.375
$ws.Range("H76").Value = 3889
$ws.Range("I76").Value = 3400
$ws.Range("J76").Value = 3943.3333
$ws.Range("K76").Value = 3400
$ws.Range("L76").Value = 3943.3333
$ws.Range("M76").Value = -3085
$ws.Range("N76").Value = -4573.3333
$ws.Range("H77").Value = 3885.8333
$ws.Range("I77").Value = 2377.5
$ws.Range("J77").Value = 4074.375
$ws.Range("K77").Value = 11887.5
$ws.Range("L77").Value = 20371.875
$ws.Range("M77").Value = -7207.5
$ws.Range("N77").Value = -29731.875
$ws.Range("H79").Value = 3889
$ws.Range("I79").Value = 3400
$ws.Range("J79").Value = 3943.3333
$ws.Range("K79").Value = 3400
$ws.Range("L79").Value = 3943.3333
$ws.Range("M79").Value = -2308
$ws.Range("N79").Value = -6127.3333
$ws.Range("H88").Value = 1179.5238
$ws.Range("I88").Value = 796
$ws.Range("K88").Value = 796
$ws.Range("M88").Value = -390
$ws.Range("H91").Value = 1179.5238
$ws.Range("I91").Value = 796
$ws.Range("K91").Value = 796
$ws.Range("M91").Value = 608
$ws.Range("H116").Value = 5958.143
$ws.Range("J116").Value = 6485.3335
$ws.Range("L116").Value = 6485.3335
$ws.Range("N116").Value = -13369.3335
$ws.Range("H132").Value = 3581.6155
$ws.Range("I132").Value = 3815.0476
$ws.Range("J132").Value = 2601.2
$ws.Range("K132").Value = 11445.1428
$ws.Range("L132").Value = 7803.599999999999
$ws.Range("M132").Value = -8915.1428
$ws.Range("N132").Value = -12863.6
$ws.Range("H137").Value = 87000
$ws.Range("I137").Value = 5349.8335
$ws.Range("K137").Value = 16049.5005
$ws.Range("M137").Value = -13499.5005

# --- Sheet: ARM (33 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1184.3438
$ws.Range("I2").Value = 1126.6364
$ws.Range("J2").Value = 1311.3
$ws.Range("K2").Value = 1126.6364
$ws.Range("L2").Value = 1311.3
$ws.Range("M2").Value = -1013.6364
$ws.Range("N2").Value = -1537.3
$ws.Range("H32").Value = 27437.977
$ws.Range("I32").Value = 27895.975
$ws.Range("K32").Value = 27895.975
$ws.Range("M32").Value = -27608.975
$ws.Range("H45").Value = 2017.0385
$ws.Range("I45").Value = 2075.2666
$ws.Range("K45").Value = 2075.2666
$ws.Range("M45").Value = -1698.2666
$ws.Range("H116").Value = 1184.3438
$ws.Range("I116").Value = 1126.6364
$ws.Range("J116").Value = 1311.3
$ws.Range("K116").Value = 1126.6364
$ws.Range("L116").Value = 1311.3
$ws.Range("M116").Value = 1167.3636
$ws.Range("N116").Value = -5899.3
$ws.Range("H122").Value = 2331.8
$ws.Range("I122").Value = 2436.5881
$ws.Range("J122").Value = 1738
$ws.Range("K122").Value = 7309.7643
$ws.Range("L122").Value = 5214
$ws.Range("M122").Value = -4859.7643
$ws.Range("N122").Value = -10114
$ws.Range("H132").Value = 13126.182
$ws.Range("I132").Value = 1557.8857
$ws.Range("K132").Value = 4673.6571
$ws.Range("M132").Value = -2143.6571

# --- Sheet: BSM (14 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1184.3438
$ws.Range("I3").Value = 1126.6364
$ws.Range("J3").Value = 1311.3
$ws.Range("K3").Value = 1126.6364
$ws.Range("L3").Value = 1311.3
$ws.Range("M3").Value = -1012.6364
$ws.Range("N3").Value = -1539.3
$ws.Range("H99").Value = 2076.0908
$ws.Range("I99").Value = 1687.4
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 1687.4
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -189.4000000000001
$ws.Range("N99").Value = -5396

# --- Sheet: CRP (29 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14048.23
$ws.Range("I31").Value = 18844.354
$ws.Range("J31").Value = 4988.8887
$ws.Range("K31").Value = 18844.354
$ws.Range("L31").Value = 4988.8887
$ws.Range("M31").Value = -18549.354
$ws.Range("N31").Value = -5578.8887
$ws.Range("H34").Value = 14048.23
$ws.Range("I34").Value = 18844.354
$ws.Range("J34").Value = 4988.8887
$ws.Range("K34").Value = 18844.354
$ws.Range("L34").Value = 4988.8887
$ws.Range("M34").Value = -18642.354
$ws.Range("N34").Value = -5392.8887
$ws.Range("H58").Value = 14157.211
$ws.Range("I58").Value = 1013.6786
$ws.Range("K58").Value = 1013.6786
$ws.Range("M58").Value = -810.6786
$ws.Range("H132").Value = 14982.052
$ws.Range("I132").Value = 19842.037
$ws.Range("J132").Value = 4047.0833
$ws.Range("K132").Value = 59526.111
$ws.Range("L132").Value = 12141.2499
$ws.Range("M132").Value = -56996.111
$ws.Range("N132").Value = -17201.2499
$ws.Range("H136").Value = 14157.211
$ws.Range("I136").Value = 1013.6786
$ws.Range("K136").Value = 3041.0358
$ws.Range("M136").Value = -491.0357999999997

# --- Sheet: CUL (40 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 50751.5
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 100003
$ws.Range("K68").Value = 4500
$ws.Range("L68").Value = 300009
$ws.Range("M68").Value = -3689
$ws.Range("N68").Value = -301631
$ws.Range("H71").Value = 50751.5
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 100003
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 900027
$ws.Range("M71").Value = -9444
$ws.Range("N71").Value = -908139
$ws.Range("H106").Value = 4574.773
$ws.Range("J106").Value = 4574.773
$ws.Range("L106").Value = 13724.319
$ws.Range("N106").Value = -15616.319
$ws.Range("H123").Value = 4377.5
$ws.Range("J123").Value = 7245
$ws.Range("L123").Value = 21735
$ws.Range("N123").Value = -26635
$ws.Range("H131").Value = 715.03
$ws.Range("J131").Value = 744
$ws.Range("L131").Value = 2232
$ws.Range("N131").Value = -12312
$ws.Range("H139").Value = 2528.647
$ws.Range("I139").Value = 2057
$ws.Range("J139").Value = 2858.8
$ws.Range("K139").Value = 6171
$ws.Range("L139").Value = 8576.400000000001
$ws.Range("M139").Value = -1031
$ws.Range("N139").Value = -18856.4
$ws.Range("H140").Value = 2614.2307
$ws.Range("I140").Value = 2372.2222
$ws.Range("J140").Value = 3158.75
$ws.Range("K140").Value = 7116.6666
$ws.Range("L140").Value = 9476.25
$ws.Range("M140").Value = -1936.6666
$ws.Range("N140").Value = -19836.25

# --- Sheet: GSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 4000
$ws.Range("J18").Value = 4000
$ws.Range("L18").Value = 4000
$ws.Range("N18").Value = -4586
$ws.Range("H132").Value = 62017.58
$ws.Range("I132").Value = 54447.95
$ws.Range("J132").Value = 87249.664
$ws.Range("K132").Value = 163343.85
$ws.Range("L132").Value = 261748.992
$ws.Range("M132").Value = -160813.85
$ws.Range("N132").Value = -266808.992

# --- Sheet: LTW (43 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1068437.5
$ws.Range("I2").Value = 1178571.4
$ws.Range("J2").Value = 297500
$ws.Range("K2").Value = 1178571.4
$ws.Range("L2").Value = 297500
$ws.Range("M2").Value = -1178459.4
$ws.Range("N2").Value = -297724
$ws.Range("H5").Value = 4000
$ws.Range("J5").Value = 4000
$ws.Range("L5").Value = 4000
$ws.Range("N5").Value = -4226
$ws.Range("H7").Value = 6770.8237
$ws.Range("I7").Value = 7449.5835
$ws.Range("J7").Value = 5141.8
$ws.Range("K7").Value = 7449.5835
$ws.Range("L7").Value = 5141.8
$ws.Range("M7").Value = -7337.5835
$ws.Range("N7").Value = -5365.8
$ws.Range("H55").Value = 72.59999999999999
$ws.Range("I55").Value = 38.5
$ws.Range("J55").Value = 106.7
$ws.Range("K55").Value = 38.5
$ws.Range("L55").Value = 106.7
$ws.Range("M55").Value = 134.5
$ws.Range("N55").Value = -452.7
$ws.Range("H126").Value = 6770.8237
$ws.Range("I126").Value = 7449.5835
$ws.Range("J126").Value = 5141.8
$ws.Range("K126").Value = 22348.7505
$ws.Range("L126").Value = 15425.4
$ws.Range("M126").Value = -19878.7505
$ws.Range("N126").Value = -20365.4
$ws.Range("H132").Value = 2549.0908
$ws.Range("I132").Value = 1786.375
$ws.Range("K132").Value = 5359.125
$ws.Range("M132").Value = -2829.125
$ws.Range("H136").Value = 47682
$ws.Range("I136").Value = 101400.4
$ws.Range("J136").Value = 2916.6667
$ws.Range("K136").Value = 304201.2
$ws.Range("L136").Value = 8750.000100000001
$ws.Range("M136").Value = -301651.2
$ws.Range("N136").Value = -13850.0001

# --- Sheet: WVR (36 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4728.4443
$ws.Range("I62").Value = 3277
$ws.Range("J62").Value = 5143.143
$ws.Range("K62").Value = 3277
$ws.Range("L62").Value = 5143.143
$ws.Range("M62").Value = -2653
$ws.Range("N62").Value = -6391.143
$ws.Range("H65").Value = 4728.4443
$ws.Range("I65").Value = 3277
$ws.Range("J65").Value = 5143.143
$ws.Range("K65").Value = 16385
$ws.Range("L65").Value = 25715.715
$ws.Range("M65").Value = -13265
$ws.Range("N65").Value = -31955.715
$ws.Range("H86").Value = 23000
$ws.Range("J86").Value = 23000
$ws.Range("L86").Value = 23000
$ws.Range("N86").Value = -25246
$ws.Range("H89").Value = 23000
$ws.Range("J89").Value = 23000
$ws.Range("L89").Value = 115000
$ws.Range("N89").Value = -126232
$ws.Range("H122").Value = 1550.1143
$ws.Range("I122").Value = 1649.826
$ws.Range("J122").Value = 1359
$ws.Range("K122").Value = 4949.478
$ws.Range("L122").Value = 4077
$ws.Range("M122").Value = -2499.478
$ws.Range("N122").Value = -8977
$ws.Range("H132").Value = 1712.9333
$ws.Range("I132").Value = 969.1
$ws.Range("J132").Value = 3200.6
$ws.Range("K132").Value = 2907.3
$ws.Range("L132").Value = 9601.799999999999
$ws.Range("M132").Value = -377.3000000000002
$ws.Range("N132").Value = -14661.8
